$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Run 50" column (AZ). This shifts the old "Mean" column (BA) left into AZ.
$ws.Range("AZ:AZ").Delete()

# Rename header "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Update column A values (Gen -> MaxFES fractional progress values)
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Update the recomputed "Mean" column (now in AZ after the delete) with the
# new mean values (averaged over the 50 remaining run columns, B:AY).
$ws.Range("AZ2").Value = 15154544560.50581
$ws.Range("AZ3").Value = 14555019627.25582
$ws.Range("AZ4").Value = 6793171561.601319
$ws.Range("AZ5").Value = 212447936.6393326
$ws.Range("AZ6").Value = 19824611.29768828
$ws.Range("AZ7").Value = 2633445.51682124
$ws.Range("AZ8").Value = 459817.96302314
$ws.Range("AZ9").Value = 58375.46750599
$ws.Range("AZ10").Value = 6271.92088064
$ws.Range("AZ11").Value = 599.80713064
$ws.Range("AZ12").Value = 39.20309848
$ws.Range("AZ13").Value = 3.03261127
$ws.Range("AZ14").Value = 0.22606965
